$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.021.44"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").Value = "3.369.00"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.88"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.98"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +2.50%  "

$ws.Range("D8").Value = "3.359.51"
$ws.Range("E8").Value = "  +0.80%  "

$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  +6.59%  "

$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.63"
$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("E13").Value = "  +3.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.13"
$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("D15").Value = "3.904.31"
$ws.Range("E15").Value = "  +0.81%  "

$ws.Range("E16").Value = "  +2.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.23"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").Value = "3.353.03"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "65.121.48"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("E21").Value = "  +1.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "452.00"
$ws.Range("E22").Value = "  +3.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.93"
$ws.Range("E23").Value = "  -2.89%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.24"
$ws.Range("E25").Value = "  +3.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.64"
$ws.Range("E26").Value = "  +2.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.75"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("E28").Value = "  +1.38%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "31.20"
$ws.Range("E29").Value = "  +4.93%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.65"
$ws.Range("E30").Value = "  -0.78%  "

$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "62.82"
$ws.Range("E32").Value = "  +7.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.44"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "576.09"
$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.65"
$ws.Range("E37").Value = "  +4.39%  "

$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.371"
$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("D41").Value = "0.0₃0739"
$ws.Range("E41").Value = "  -1.25%  "

$ws.Range("D42").Value = "3.080.26"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("E44").Value = "  -0.60%  "

$ws.Range("E45").Value = "  +3.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  -0.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.16"
$ws.Range("E47").Value = "  -2.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.81"
$ws.Range("E49").Value = "  +5.17%  "

$ws.Range("E50").Value = "  -2.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.28"
$ws.Range("E51").Value = "  -0.04%  "

